$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 29, 30, 31), continuing the 2020 biweekly
# water-quality monitoring series for 청평 station.
$rows = @(
    @{
        r = 29
        A = "2020-07-15 12:00"; B = "3회차"; AB = "20200715"
        AC = 201.9; AD = 24; AE = 7.9; AF = 8.7; AG = 0.7; AH = 3.1; AI = 0.6; AJ = 127
        AK = 1.104; AL = 1.097; AM = 0.932; AN = 0.042; AO = 0.011; AP = 0.008; AQ = 0.002
        AR = 2.5; AS = 1.7
    },
    @{
        r = 30
        A = "2020-07-20 12:00"; B = "4회차"; AB = "20200722"
        AC = 121.2; AD = 25.1; AE = 7.7; AF = 9.3; AG = 0.6; AH = 3.5; AI = 1; AJ = 120
        AK = 1.057; AL = 1.04; AM = 0.899; AN = 0.032; AO = 0.011; AP = 0.008; AQ = 0.005
        AR = 4.9; AS = 1.7
    },
    @{
        r = 31
        A = "2020-07-25 12:00"; B = "5회차"; AB = "20200729"
        AC = 165.6; AD = 23.2; AE = 8.2; AF = 11.5; AG = 2.2; AH = 4.7; AI = 4.2; AJ = 114
        AK = 1.672; AL = 1.575; AM = 1.506; AN = 0.007; AO = 0.026; AP = 0.013; AQ = 0.001
        AR = 31.2; AS = 2.4
    }
)

foreach ($row in $rows) {
    $r = $row.r
    $src = $r - 1

    # Copy the whole previous row's formatting/formulas as a template, then
    # overwrite with this row's specific values.
    $ws.Range("A$src`:BJ$src").Copy()
    $ws.Range("A$r`:BJ$r").PasteSpecial(-4104)

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 28).Value = $row.AB

    $ws.Cells.Item($r, 29).Value = $row.AC
    $ws.Cells.Item($r, 30).Value = $row.AD
    $ws.Cells.Item($r, 31).Value = $row.AE
    $ws.Cells.Item($r, 32).Value = $row.AF
    $ws.Cells.Item($r, 33).Value = $row.AG
    $ws.Cells.Item($r, 34).Value = $row.AH
    $ws.Cells.Item($r, 35).Value = $row.AI
    $ws.Cells.Item($r, 36).Value = $row.AJ
    $ws.Cells.Item($r, 37).Value = $row.AK
    $ws.Cells.Item($r, 38).Value = $row.AL
    $ws.Cells.Item($r, 39).Value = $row.AM
    $ws.Cells.Item($r, 40).Value = $row.AN
    $ws.Cells.Item($r, 41).Value = $row.AO
    $ws.Cells.Item($r, 42).Value = $row.AP
    $ws.Cells.Item($r, 43).Value = $row.AQ
    $ws.Cells.Item($r, 44).Value = $row.AR
    $ws.Cells.Item($r, 45).Value = $row.AS
}

$ws.Range("F30").Select()
